$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.06515966666666667
$ws.Range("M2").Value = 2.33201
$ws.Range("N2").Value = 6.99603
$ws.Range("O2").Value = 0.3303263034789547
$ws.Range("P2").Value = 0.3303263034789548
$ws.Range("Q2").Value = 0.1519529942633333
$ws.Range("R2").Value = 1.36757694837
$ws.Range("S2").Value = 0.3303263034789547
$ws.Range("T2").Value = 0.3303263034789548

# Row 3
$ws.Range("G3").Value = 0.06515966666666667
$ws.Range("N3").Value = 5.238131999999999
$ws.Range("O3").Value = 0.2473249515360603
$ws.Range("P3").Value = 0.2473249515360603
$ws.Range("S3").Value = 0.2473249515360603
$ws.Range("T3").Value = 0.2473249515360603

# Row 4
$ws.Range("G4").Value = 0.06515966666666667
$ws.Range("M4").Value = 1.145780666666667
$ws.Range("N4").Value = 3.437342
$ws.Range("O4").Value = 0.1622984001859565
$ws.Range("P4").Value = 0.1622984001859565
$ws.Range("Q4").Value = 0.07465868631311111
$ws.Range("R4").Value = 0.671928176818
$ws.Range("S4").Value = 0.1622984001859565
$ws.Range("T4").Value = 0.1622984001859565

# Row 5
$ws.Range("G5").Value = 0.06515966666666667
$ws.Range("M5").Value = 1.835881666666667
$ws.Range("N5").Value = 5.507645
$ws.Range("O5").Value = 0.2600503447990285
$ws.Range("P5").Value = 0.2600503447990286
$ws.Range("Q5").Value = 0.1196254374394445
$ws.Range("R5").Value = 1.076628936955
$ws.Range("S5").Value = 0.2600503447990285
$ws.Range("T5").Value = 0.2600503447990286
